$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 9024922
$ws.Range("I80").Value = 464.16666
$ws.Range("J80").Value = 16244489
$ws.Range("K80").Value = 1392.49998
$ws.Range("L80").Value = 48733467
$ws.Range("M80").Value = -394.4999800000001
$ws.Range("N80").Value = -48735463
$ws.Range("H83").Value = 9024922
$ws.Range("I83").Value = 464.16666
$ws.Range("J83").Value = 16244489
$ws.Range("K83").Value = 4177.49994
$ws.Range("L83").Value = 146200401
$ws.Range("M83").Value = 814.5000600000003
$ws.Range("N83").Value = -146210385
$ws.Range("H94").Value = 3000
$ws.Range("I94").Value = 3000
$ws.Range("K94").Value = 3000
$ws.Range("M94").Value = -2549
$ws.Range("H129").Value = 323535.2
$ws.Range("J129").Value = 358171.22
$ws.Range("L129").Value = 1074513.66
$ws.Range("N129").Value = -1084513.66
$ws.Range("H138").Value = 3874.1562
$ws.Range("I138").Value = 2732.4614
$ws.Range("J138").Value = 4655.316
$ws.Range("K138").Value = 8197.3842
$ws.Range("L138").Value = 13965.948
$ws.Range("M138").Value = -3057.3842
$ws.Range("N138").Value = -24245.948
$ws.Range("H141").Value = 1685.1522
$ws.Range("I141").Value = 1542.925
$ws.Range("J141").Value = 2633.3333
$ws.Range("K141").Value = 4628.775
$ws.Range("L141").Value = 7899.999899999999
$ws.Range("M141").Value = 551.2250000000004
$ws.Range("N141").Value = -18259.9999
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H6").Value = 6253375
$ws.Range("I6").Value = 12501750
$ws.Range("K6").Value = 12501750
$ws.Range("M6").Value = -12501577
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 0
$ws.Range("K25").Value = 0
$ws.Range("M25").Value = ""
$ws.Range("H32").Value = 8936.406000000001
$ws.Range("I32").Value = 6836.4683
$ws.Range("J32").Value = 18694.941
$ws.Range("K32").Value = 6836.4683
$ws.Range("L32").Value = 18694.941
$ws.Range("M32").Value = -6549.4683
$ws.Range("N32").Value = -19268.941
$ws.Range("H39").Value = 3000
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = 3000
$ws.Range("K39").Value = 0
$ws.Range("L39").Value = 3000
$ws.Range("M39").Value = ""
$ws.Range("N39").Value = -4040
$ws.Range("H61").Value = 9262528
$ws.Range("I61").Value = 11114268
$ws.Range("J61").Value = 3831.1667
$ws.Range("K61").Value = 11114268
$ws.Range("L61").Value = 3831.1667
$ws.Range("M61").Value = -11114056
$ws.Range("N61").Value = -4255.1667
$ws.Range("H63").Value = 3907448.5
$ws.Range("I63").Value = 1357.6
$ws.Range("K63").Value = 1357.6
$ws.Range("M63").Value = -671.5999999999999
$ws.Range("H66").Value = 3907448.5
$ws.Range("I66").Value = 1357.6
$ws.Range("K66").Value = 6788
$ws.Range("M66").Value = -3356
$ws.Range("H74").Value = 27779292
$ws.Range("I74").Value = 41667376
$ws.Range("J74").Value = 3124.1667
$ws.Range("K74").Value = 41667376
$ws.Range("L74").Value = 3124.1667
$ws.Range("M74").Value = -41666502
$ws.Range("N74").Value = -4872.1667
$ws.Range("H77").Value = 27779292
$ws.Range("I77").Value = 41667376
$ws.Range("J77").Value = 3124.1667
$ws.Range("K77").Value = 208336880
$ws.Range("L77").Value = 15620.8335
$ws.Range("M77").Value = -208332512
$ws.Range("N77").Value = -24356.8335
$ws.Range("H131").Value = 60715
$ws.Range("J131").Value = 60715
$ws.Range("L131").Value = 60715
$ws.Range("N131").Value = -70795
$ws.Range("H132").Value = 12515088
$ws.Range("I132").Value = 16131379
$ws.Range("J132").Value = 58972.668
$ws.Range("K132").Value = 48394137
$ws.Range("L132").Value = 176918.004
$ws.Range("M132").Value = -48391607
$ws.Range("N132").Value = -181978.004
$ws.Range("H136").Value = 9262528
$ws.Range("I136").Value = 11114268
$ws.Range("J136").Value = 3831.1667
$ws.Range("K136").Value = 33342804
$ws.Range("L136").Value = 11493.5001
$ws.Range("M136").Value = -33340254
$ws.Range("N136").Value = -16593.5001
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3964.182
$ws.Range("I20").Value = 3450.75
$ws.Range("K20").Value = 3450.75
$ws.Range("M20").Value = -3203.75
$ws.Range("H81").Value = 14684.375
$ws.Range("J81").Value = 14684.375
$ws.Range("L81").Value = 14684.375
$ws.Range("N81").Value = -16806.375
$ws.Range("H84").Value = 14684.375
$ws.Range("J84").Value = 14684.375
$ws.Range("L84").Value = 44053.125
$ws.Range("N84").Value = -54661.125
$ws.Range("H86").Value = 2747.9333
$ws.Range("I86").Value = 2439.3845
$ws.Range("J86").Value = 4753.5
$ws.Range("K86").Value = 2439.3845
$ws.Range("L86").Value = 4753.5
$ws.Range("M86").Value = -1316.3845
$ws.Range("N86").Value = -6999.5
$ws.Range("H89").Value = 2747.9333
$ws.Range("I89").Value = 2439.3845
$ws.Range("J89").Value = 4753.5
$ws.Range("K89").Value = 12196.9225
$ws.Range("L89").Value = 23767.5
$ws.Range("M89").Value = -6580.922500000001
$ws.Range("N89").Value = -34999.5
$ws.Range("H94").Value = 943.85187
$ws.Range("I94").Value = 480.94116
$ws.Range("K94").Value = 480.94116
$ws.Range("M94").Value = -29.94116000000002
$ws.Range("H99").Value = 1262.7273
$ws.Range("J99").Value = 800
$ws.Range("L99").Value = 800
$ws.Range("N99").Value = -3796
$ws.Range("H105").Value = 1284104.4
$ws.Range("I105").Value = 1436.1538
$ws.Range("J105").Value = 1925438.4
$ws.Range("K105").Value = 1436.1538
$ws.Range("L105").Value = 1925438.4
$ws.Range("M105").Value = 310.8462
$ws.Range("N105").Value = -1928932.4
$ws.Range("H110").Value = 22762.25
$ws.Range("I110").Value = 10000
$ws.Range("J110").Value = 35524.5
$ws.Range("K110").Value = 10000
$ws.Range("L110").Value = 35524.5
$ws.Range("M110").Value = -5910
$ws.Range("N110").Value = -43704.5
$ws.Range("H130").Value = 43779.668
$ws.Range("J130").Value = 43779.668
$ws.Range("L130").Value = 43779.668
$ws.Range("N130").Value = -53819.668
$ws.Range("H135").Value = 30619.25
$ws.Range("J135").Value = 30619.25
$ws.Range("L135").Value = 30619.25
$ws.Range("N135").Value = -40759.25
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 147.09091
$ws.Range("J22").Value = 175
$ws.Range("L22").Value = 175
$ws.Range("N22").Value = -875
$ws.Range("H132").Value = 31252450
$ws.Range("I132").Value = 38463080
$ws.Range("J132").Value = 6385.3335
$ws.Range("K132").Value = 115389240
$ws.Range("L132").Value = 19156.0005
$ws.Range("M132").Value = -115386710
$ws.Range("N132").Value = -24216.0005
$ws.Range("H141").Value = 22364
$ws.Range("J141").Value = 22364
$ws.Range("L141").Value = 22364
$ws.Range("N141").Value = -32724
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H15").Value = 69.59999999999999
$ws.Range("I15").Value = 50
$ws.Range("J15").Value = 99
$ws.Range("K15").Value = 150
$ws.Range("L15").Value = 297
$ws.Range("M15").Value = -10
$ws.Range("N15").Value = -577
$ws.Range("H129").Value = 223458.22
$ws.Range("J129").Value = 401560
$ws.Range("L129").Value = 1204680
$ws.Range("N129").Value = -1214680
$ws.Range("H131").Value = 708.77
$ws.Range("J131").Value = 752.5909
$ws.Range("L131").Value = 2257.7727
$ws.Range("N131").Value = -12337.7727
$ws.Range("H132").Value = 1000
$ws.Range("J132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("N132").Value = ""
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 3680891.2
$ws.Range("I70").Value = 4699.857
$ws.Range("K70").Value = 4699.857
$ws.Range("M70").Value = -4429.857
$ws.Range("H73").Value = 3680891.2
$ws.Range("I73").Value = 4699.857
$ws.Range("K73").Value = 4699.857
$ws.Range("M73").Value = -3763.857
$ws.Range("H132").Value = 5104475
$ws.Range("I132").Value = 6688983
$ws.Range("J132").Value = 86866.164
$ws.Range("K132").Value = 20066949
$ws.Range("L132").Value = 260598.492
$ws.Range("M132").Value = -20064419
$ws.Range("N132").Value = -265658.492
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 1944.6666
$ws.Range("I93").Value = 2014.1666
$ws.Range("K93").Value = 2014.1666
$ws.Range("M93").Value = -766.1666
$ws.Range("H100").Value = 2254.3845
$ws.Range("I100").Value = 1961.4
$ws.Range("K100").Value = 1961.4
$ws.Range("M100").Value = -1420.4
$ws.Range("H110").Value = 40009
$ws.Range("J110").Value = 40009
$ws.Range("L110").Value = 40009
$ws.Range("N110").Value = -48189
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1857.04
$ws.Range("I126").Value = 1416.8948
$ws.Range("J126").Value = 3250.8333
$ws.Range("K126").Value = 4250.6844
$ws.Range("L126").Value = 9752.499899999999
$ws.Range("M126").Value = -1780.6844
$ws.Range("N126").Value = -14692.4999
